# Finalized "Softwares disponiveis" INSERT sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INSERT")

# The sheet used to have a gap of empty rows (2-5) between the A1 prefix
# cell and the header row (old row 6). Remove that gap so the header
# becomes row 2 and the data starts at row 3 (shifting everything up by 4).
$ws.Rows("2:5").Delete()

# A1 holds the literal SQL INSERT prefix - it now needs two more columns.
$ws.Range("A1").Value = "INSERT INTO softwares_disponiveis (nome_software, descricao_software, nome_arquivo, nome_imagem) values("

# Column H used to hold the old combined formula (local-only version); it is
# no longer part of the table, so clear it out entirely.
$ws.Range("H1:H24").Clear()

# New header row (row 2): SOFTWARE | DESCRICAO | CODIGO_SQL | NOME_ARQUIVO | NOME_IMAGEM
$ws.Range("D2").Value = "CÓDIGO_SQL"
$ws.Range("E2").Value = "NOME_ARQUIVO"
$ws.Range("F2").Value = "NOME_IMAGEM"

$q = [char]34

# Per-row file name / image name / status data.
$rowsData = @(
    @{ Row=3;  Arquivo="android-studio-bundle-135.1641136.exe";    Imagem="android-studio.jpg";  Status="ok" },
    @{ Row=4;  Arquivo="astah-community-6_8_0-d254c5-jre-setup";   Imagem="astah.png";            Status="ok" },
    @{ Row=5;  Arquivo="CS2_setup.exe";                            Imagem="case_studio.gif";      Status="ok" },
    @{ Row=6;  Arquivo="packettracer533";                          Imagem="packet_tracer.jpg";    Status="ok" },
    @{ Row=7;  Arquivo="Dev-Cpp 64 bits";                          Imagem="dev_c++.jpg";          Status="ok" },
    @{ Row=8;  Arquivo="dia-setup-0.97.2-2.exe";                   Imagem="dia.jpg";              Status="-"  },
    @{ Row=9;  Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=10; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=11; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=12; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=13; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=14; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=15; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=16; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=17; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=18; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=19; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=20; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=21; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=22; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=23; Arquivo="-";                                        Imagem="-";                    Status="-"  },
    @{ Row=24; Arquivo="-";                                        Imagem="-";                    Status="-"  }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Range("D$r").Value = $rd.Arquivo
    $ws.Range("E$r").Value = $rd.Imagem

    $fargs = @("`$A`$1", "$q'$q", "B$r", "$q'$q", "$q, $q", "$q'$q", "C$r", "$q'$q", "$q, $q", " $q'$q", "D$r", "$q'$q", "$q, $q", " $q'$q", "E$r", "$q'$q", "$q);$q")
    $formula = "=CONCATENATE(" + ($fargs -join ",") + ")"
    $ws.Range("F$r").Formula = $formula

    $ws.Range("G$r").Value = $rd.Status
}

# Conditional formatting: highlight "ok" cells in column G (rows 3-24) with
# the built-in "light red fill / dark red text" style.
$rng = $ws.Range("G3:G24")
$fc = $rng.FormatConditions.Add(9, 0, "ok", "ok")
$fc.Text = "ok"
$fc.Formula1 = "=NOT(ISERROR(SEARCH(""ok"",G3)))"
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# Final selection/cursor position left on the sheet.
$ws.Range("F8").Select()
